$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 ("Why is Mobile Device Security Important?"):
# add a small source-link textbox at the bottom of the slide.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tb1 = $s3.Shapes.AddTextbox(1, -0.0015748031496062992, 519.8917322834645, 960.0, 19.38748031496063)
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1
$tb1.Fill.Visible = $false

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "https://"
$tr1.InsertAfter("antivirus.comodo.com") | Out-Null
$tr1.InsertAfter("/blog/comodo-news/why-mobile-security-is-important-today/") | Out-Null
$tr1.Font.Size = 10
$tb1.Height = 19.3875

# ---------------------------------------------------------------------------
# Slide 4 ("Dangers"):
# append two more bullet paragraphs to the content placeholder, then add a
# source-link textbox at the bottom of the slide.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$content4 = $s4.Shapes.Item(4)
$content4.TextFrame.TextRange.InsertAfter("`rSpyware`rWeak Passwords") | Out-Null

$tb2 = $s4.Shapes.AddTextbox(1, 0.0, 519.8916535433071, 952.8, 19.38748031496063)
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1
$tb2.Fill.Visible = $false

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "https://auth0.com/blog/the-9-most-common-security-threats-to-mobile-devices-in-2021/"
$tr2.Font.Size = 10
$tb2.Height = 19.3875

# ---------------------------------------------------------------------------
# Slide 5 ("Best Practices"):
# add a source-link textbox at the bottom of the slide.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tb3 = $s5.Shapes.AddTextbox(1, 0.0, 519.8917322834645, 952.8, 19.38748031496063)
$tb3.TextFrame.WordWrap = $true
$tb3.TextFrame.AutoSize = 1
$tb3.Fill.Visible = $false

$tr3 = $tb3.TextFrame.TextRange
$tr3.Text = "https://"
$tr3.InsertAfter("www.ntiva.com") | Out-Null
$tr3.InsertAfter("/blog/top-7-mobile-device-security-best-practices") | Out-Null
$tr3.Font.Size = 10
$tb3.Height = 19.3875
